# Update the "want-to-go" counts (想去人数, column F) on the sheets
# that list exhibitions/activities: "展览" and "全部类型".
#
# Row (by name in col C)                         old -> new
#  南宁·熊喵M动漫嘉年华【免费】        (row 2)    1418 -> 1419
#  南宁·第二届北极光动漫展              (row 3)    2981 -> 2983
#  南宁·原神x星铁x绝区零同人ONLY3.0     (row 4)      28 ->   29
#  南宁·万圣漫控嘉年华10                              276 ->  278

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # "展览" sheet has 万圣漫控嘉年华10 on row 5, "全部类型" has it on row 6,
    # so locate each row by matching the event name in column C instead of
    # hard-coding row numbers.
    $updates = @{
        "南宁·熊喵M动漫嘉年华【免费】" = 1419
        "南宁·第二届北极光动漫展" = 2983
        "南宁·原神x星铁x绝区零同人ONLY3.0" = 29
        "南宁·万圣漫控嘉年华10" = 278
    }

    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count

    for ($r = 2; $r -le $rowCount; $r++) {
        $eventName = $ws.Cells.Item($r, 3).Value2
        if ($updates.ContainsKey($eventName)) {
            $ws.Cells.Item($r, 6).Value = $updates[$eventName]
        }
    }
}
